# Generate Report for Handoff
#
# "b.md" has been handed off again (a new xlf was generated from a newer
# source commit), so its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff" across the
# Overview / zh-cn / de-de sheets. Each locale sheet also records the new
# handoff file name/time, flips "Content Duplicate" to False, widens the
# "Error Detail" column and fills in a stale-handback warning there.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d4d1bb80ac946272d30321f8783690af0491d389/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e98f4a674711513594e649fbb8216cd2af824652/e2e/b.md."

# ---- Overview sheet: row 3 is the "b.md" file ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-28 14:37:25"

# ---- zh-cn sheet: row 3 is the "b.md" file ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# "True"/"False" look like booleans, but the workbook stores them as plain
# text in the shared-string table (see F2, which already holds text
# "False"). Copy that existing text cell instead of assigning the literal
# string, so the destination keeps the text type instead of being
# auto-coerced to a real Boolean.
$zhcn.Range("F2").Copy()
$zhcn.Range("F3").PasteSpecial(-4163)
$zhcn.Application.CutCopyMode = $false
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-28 14:37:21"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.1666667

# ---- de-de sheet: row 3 is the "b.md" file ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F2").Copy()
$dede.Range("F3").PasteSpecial(-4163)
$dede.Application.CutCopyMode = $false
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-28 14:37:25"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.1666667
